$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the previously-empty predicted price cell (C2) with a formula
# that derives hours from the minutes value already in B2.
$ws.Range("C2").Formula = "=B2/60"
